$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46006.01041666666, 937.261),
    @(3, 46006.02083333334, 932.956),
    @(4, 46006.03125, 929.153),
    @(5, 46006.04166666666, 925.774),
    @(6, 46006.05208333334, 913.081),
    @(7, 46006.0625, 909.109),
    @(8, 46006.07291666666, 905.133),
    @(9, 46006.08333333334, 900.58),
    @(10, 46006.09375, 884.306),
    @(11, 46006.10416666666, 876.432),
    @(12, 46006.11458333334, 869.972),
    @(13, 46006.125, 863.168),
    @(14, 46006.13541666666, 837.7),
    @(15, 46006.14583333334, 826.774),
    @(16, 46006.15625, 816.167),
    @(17, 46006.16666666666, 806.5700000000001),
    @(18, 46006.17708333334, 750.4589999999999),
    @(19, 46006.1875, 736.3630000000001),
    @(20, 46006.19791666666, 723.223),
    @(21, 46006.20833333334, 711.098),
    @(22, 46006.21875, 632.876),
    @(23, 46006.22916666666, 622.7329999999999),
    @(24, 46006.23958333334, 610.134),
    @(25, 46006.25, 597.129),
    @(26, 46006.26041666666, 524.886),
    @(27, 46006.27083333334, 513.283),
    @(28, 46006.28125, 503.74),
    @(29, 46006.29166666666, 493.923),
    @(30, 46006.30208333334, 436.062),
    @(31, 46006.3125, 421.035),
    @(32, 46006.32291666666, 411.93),
    @(33, 46006.33333333334, 402.774),
    @(34, 46006.34375, 351.746),
    @(35, 46006.35416666666, 346.643),
    @(36, 46006.36458333334, 341.089),
    @(37, 46006.375, 334.89),
    @(38, 46006.38541666666, 306.955),
    @(39, 46006.39583333334, 300.059),
    @(40, 46006.40625, 293.136),
    @(41, 46006.41666666666, 287.493),
    @(42, 46006.42708333334, 277.635),
    @(43, 46006.4375, 280.59),
    @(44, 46006.44791666666, 283.653),
    @(45, 46006.45833333334, 286.745),
    @(46, 46006.46875, 307.512),
    @(47, 46006.47916666666, 313.838),
    @(48, 46006.48958333334, 320.726),
    @(49, 46006.5, 326.03),
    @(50, 46006.51041666666, 342.461),
    @(51, 46006.52083333334, 342.234),
    @(52, 46006.53125, 341.676),
    @(53, 46006.54166666666, 341.601),
    @(54, 46006.55208333334, 334.062),
    @(55, 46006.5625, 331.693),
    @(56, 46006.57291666666, 328.983),
    @(57, 46006.58333333334, 326.873),
    @(58, 46006.59375, 325.235),
    @(59, 46006.60416666666, 324.801),
    @(60, 46006.61458333334, 324.689),
    @(61, 46006.625, 324.475),
    @(62, 46006.63541666666, 350.469),
    @(63, 46006.64583333334, 358.493),
    @(64, 46006.65625, 371.012),
    @(65, 46006.66666666666, 382.068),
    @(66, 46006.67708333334, 419.286),
    @(67, 46006.6875, 430.392),
    @(68, 46006.69791666666, 441.486),
    @(69, 46006.70833333334, 452.66),
    @(70, 46006.71875, 479.221),
    @(71, 46006.72916666666, 486.026),
    @(72, 46006.73958333334, 491.678),
    @(73, 46006.75, 496.765),
    @(74, 46006.76041666666, 517.583),
    @(75, 46006.77083333334, 515.527),
    @(76, 46006.78125, 516.807),
    @(77, 46006.79166666666, 518.3150000000001),
    @(78, 46006.80208333334, 522.943),
    @(79, 46006.8125, 518.777),
    @(80, 46006.82291666666, 514.928),
    @(81, 46006.83333333334, 511.411),
    @(82, 46006.84375, 502.445),
    @(83, 46006.85416666666, 496.982),
    @(84, 46006.86458333334, 492.478),
    @(85, 46006.875, 487.799),
    @(86, 46006.88541666666, 461.875),
    @(87, 46006.89583333334, 455.647),
    @(88, 46006.90625, 449.685),
    @(89, 46006.91666666666, 443.43),
    @(90, 46006.92708333334, 412.369),
    @(91, 46006.9375, 405.978),
    @(92, 46006.94791666666, 402.926),
    @(93, 46006.95833333334, 396.87),
    @(94, 46006.96875, 0),
    @(95, 46006.97916666666, 0),
    @(96, 46006.98958333334, 0),
    @(97, 46007.0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
}

